$x = 1 + 2
Write-Host "sum:" $x
for ($i = 0; $i -lt 3; $i++) {
  Write-Host "i=" $i
}
$arr = @(1,2,3)
Write-Host "arr:" $arr.Length
